$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Update client code cell (name/email) to the actual client's name and email
$ws.Range("A31").Value = "Charlie, charlie@mail.com"

# Update the "Client discount" unit price (amount column recalculates via formula)
$ws.Range("E18").Value = 125

$wb.Save()
